$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" data row (old row 26) entirely - its data is dropped,
# and all subsequent rows shift up by one.
$ws.Rows(26).Delete()

# Remove the "SC 92" data row (old row 28, now row 27 after the first
# deletion) entirely - its data is dropped too, shifting remaining rows up.
$ws.Rows(27).Delete()

# Column F ("F") holds recomputed error values for the surviving rows;
# refresh them to match the newly recalculated set.
$ws.Range("F26").Value = 17.38
$ws.Range("F27").Value = ""
$ws.Range("F28").Value = ""
$ws.Range("F29").Value = 18.06
$ws.Range("F30").Value = 16.89
$ws.Range("F31").Value = ""
$ws.Range("F32").Value = ""
$ws.Range("F33").Value = 17.53
